# Add a second worksheet ("Sheet2") after the existing "Sheet1" and populate
# it with a small Grade / Subjects / Contents table, matching the target
# workbook produced by the author's upload.

$wb = $excel.ActiveWorkbook

# Create the new sheet right after the last existing sheet (so it lands
# after Sheet1, not before it) and name it "Sheet2".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Sheet2"

# Fill in the table. The order in which new string values are first written
# determines the shared-string table order, so we deliberately write the
# "Contents" column's row-2 value (Quiz) before the "Contents" header cell,
# and Flashcards (C4) before Science (B4), to reproduce the exact shared
# string ordering of the target file.
$ws2.Range("A1").Value = "Grade"
$ws2.Range("B1").Value = "Subjects"
$ws2.Range("C2").Value = "Quiz"
$ws2.Range("C1").Value = "Contents"
$ws2.Range("B2").Value = "English"
$ws2.Range("B3").Value = "Math"
$ws2.Range("C3").Value = "Quiz,Worksheet"
$ws2.Range("C4").Value = "Flashcards"
$ws2.Range("B4").Value = "Science"

$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2

# Match the column widths recorded for the new sheet.
$ws2.Columns.Item(3).ColumnWidth = 24.44140625
$ws2.Columns.Item(4).ColumnWidth = 21.33203125

# Leave the cursor / selection on Sheet2 at K19, and make Sheet2 the active
# (selected) tab, same as the target workbook (activeTab points at Sheet2
# and tabSelected moves off of Sheet1 onto Sheet2).
$ws2.Range("K19").Select()
$ws2.Activate()
